$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy header formatting (bold, border, centered) from E1,
# then set the new header text "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F6 - plain text values (microsecond-precision timestamps
# stored as text, same as the rest of the sheet's data cells)
$ws.Range("F2").Value = "2021-10-05 13:41:56.084029"
$ws.Range("F3").Value = "2021-10-05 13:41:56.084041"
$ws.Range("F4").Value = "2021-10-05 13:41:56.084045"
$ws.Range("F5").Value = "2021-10-05 13:41:56.084048"
$ws.Range("F6").Value = "2021-10-05 13:41:56.084051"
